# Generate Report for Archive
#
# The localization run moved from "Ready for handoff" to "In Translation".
# That status string is shared by:
#   - Overview sheet, row 2, columns E (zh-cn) and F (de-de)
#   - zh-cn sheet,   row 2, column C (Status)
#   - de-de sheet,   row 2, column C (Status)
# Updating the text also shrinks those status columns (the new text is
# shorter), so their column widths are tightened to match.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update the status cells (shared string text change).
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Tighten the now-narrower status columns.
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
